# Apply the edit described by the diff:
# - Insert a new column D ("MAE") before the existing "Tipo" column (which shifts to E)
# - Update the MSE (B) and R2 (C) values with the new re-computed results
# - Fill the new MAE (D) column with its computed values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column D ("Tipo"), shifting it to column E
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("D1").Value = "MAE"

# Carry over the header style (bold, bordered, centered) used by the other header cells
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1
$ws.Range("D1").Borders.Weight = 2

# Updated MSE (B) and R2 (C) values, plus the new MAE (D) values
$ws.Range("B2").Value = 0.4990795993067688
$ws.Range("C2").Value = 0.9900622653217547
$ws.Range("D2").Value = 0.5739053831682486

$ws.Range("B3").Value = 0.2405002146899179
$ws.Range("C3").Value = 0.9952994172557589
$ws.Range("D3").Value = 0.381236035918039

$ws.Range("B4").Value = 0.3485103678959032
$ws.Range("C4").Value = 0.9932954354838006
$ws.Range("D4").Value = 0.4721310413590081

$ws.Range("B5").Value = 0.4037020542686865
$ws.Range("C5").Value = 0.9920395482072137
$ws.Range("D5").Value = 0.4999052006732372
